$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New rows 36-42 with identifiers and values (per diff).
# Column A (the new unique identifier strings) is populated first, in
# row order, so the shared-strings table gains NUA2, NUA1, NUA3, PA2,
# PA3, PB1, PC1 in that sequence; the lone text value in column B
# ("399285450#0") is written afterwards so it lands last in the table.
$ws.Cells.Item(36, 1).Value = "NUA2"
$ws.Cells.Item(37, 1).Value = "NUA1"
$ws.Cells.Item(38, 1).Value = "NUA3"
$ws.Cells.Item(39, 1).Value = "PA2"
$ws.Cells.Item(40, 1).Value = "PA3"
$ws.Cells.Item(41, 1).Value = "PB1"
$ws.Cells.Item(42, 1).Value = "PC1"

$ws.Cells.Item(36, 2).Value = 552654210
$ws.Cells.Item(37, 2).Value = "399285450#0"
$ws.Cells.Item(38, 2).Value = 48202933
$ws.Cells.Item(39, 2).Value = 305214685
$ws.Cells.Item(39, 3).Value = 542435213
$ws.Cells.Item(40, 2).Value = 542440756
$ws.Cells.Item(41, 2).Value = 542435213
$ws.Cells.Item(42, 2).Value = 542440756

# Update view/selection to match final state: scrolled so row 13 is the
# top visible row, with E38 as the active selected cell.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E38").Select()
